$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (was "Sheet1")
$ws.Name = "Yashasvi Jaiswal"

# Insert a new first column; existing A..L data shifts to B..M
$ws.Columns.Item(1).Insert()

# Pre-format the full target range as Text so numeric-looking values
# (matchNo, runs, balls, fours, sixes, sr) are stored as literal strings,
# matching the source data (all cells are t="str" in the target file).
$target = $ws.Range("A1:M11")
$target.NumberFormat = "@"

# Row 1
$ws.Cells.Item(1, 1).Value = "matchNo"
$ws.Cells.Item(1, 2).Value = "teamName"
$ws.Cells.Item(1, 3).Value = "batterName"
$ws.Cells.Item(1, 4).Value = "states"
$ws.Cells.Item(1, 5).Value = "runs"
$ws.Cells.Item(1, 6).Value = "balls"
$ws.Cells.Item(1, 7).Value = "fours"
$ws.Cells.Item(1, 8).Value = "sixes"
$ws.Cells.Item(1, 9).Value = "sr"
$ws.Cells.Item(1, 10).Value = "opponentTeamName"
$ws.Cells.Item(1, 11).Value = "venue"
$ws.Cells.Item(1, 12).Value = "date"
$ws.Cells.Item(1, 13).Value = "result"

# Row 2
$ws.Cells.Item(2, 1).Value = "28th"
$ws.Cells.Item(2, 2).Value = "Rajasthan Royals"
$ws.Cells.Item(2, 3).Value = "Yashasvi Jaiswal"
$ws.Cells.Item(2, 4).Value = "lbw b Rashid Khan"
$ws.Cells.Item(2, 5).Value = "12"
$ws.Cells.Item(2, 6).Value = "13"
$ws.Cells.Item(2, 7).Value = "2"
$ws.Cells.Item(2, 8).Value = "0"
$ws.Cells.Item(2, 9).Value = "92.30"
$ws.Cells.Item(2, 10).Value = "Sunrisers Hyderabad"
$ws.Cells.Item(2, 11).Value = "Delhi"
$ws.Cells.Item(2, 12).Value = "May 02"
$ws.Cells.Item(2, 13).Value = "Royals won by 55 runs"

# Row 3
$ws.Cells.Item(3, 1).Value = "47th"
$ws.Cells.Item(3, 2).Value = "Rajasthan Royals"
$ws.Cells.Item(3, 3).Value = "Yashasvi Jaiswal"
$ws.Cells.Item(3, 4).Value = "c †Dhoni b Asif"
$ws.Cells.Item(3, 5).Value = "50"
$ws.Cells.Item(3, 6).Value = "21"
$ws.Cells.Item(3, 7).Value = "6"
$ws.Cells.Item(3, 8).Value = "3"
$ws.Cells.Item(3, 9).Value = "238.09"
$ws.Cells.Item(3, 10).Value = "Chennai Super Kings"
$ws.Cells.Item(3, 11).Value = "Abu Dhabi"
$ws.Cells.Item(3, 12).Value = "October 02"
$ws.Cells.Item(3, 13).Value = "Royals won by 7 wickets (with 15 balls remaining)"

# Row 4
$ws.Cells.Item(4, 1).Value = "24th"
$ws.Cells.Item(4, 2).Value = "Rajasthan Royals"
$ws.Cells.Item(4, 3).Value = "Yashasvi Jaiswal"
$ws.Cells.Item(4, 4).Value = "c & b Chahar"
$ws.Cells.Item(4, 5).Value = "32"
$ws.Cells.Item(4, 6).Value = "20"
$ws.Cells.Item(4, 7).Value = "2"
$ws.Cells.Item(4, 8).Value = "2"
$ws.Cells.Item(4, 9).Value = "160.00"
$ws.Cells.Item(4, 10).Value = "Mumbai Indians"
$ws.Cells.Item(4, 11).Value = "Delhi"
$ws.Cells.Item(4, 12).Value = "April 29"
$ws.Cells.Item(4, 13).Value = "Mumbai won by 7 wickets (with 9 balls remaining)"

# Row 5
$ws.Cells.Item(5, 1).Value = "51st"
$ws.Cells.Item(5, 2).Value = "Rajasthan Royals"
$ws.Cells.Item(5, 3).Value = "Yashasvi Jaiswal"
$ws.Cells.Item(5, 4).Value = "c †Ishan Kishan b Coulter-Nile"
$ws.Cells.Item(5, 5).Value = "12"
$ws.Cells.Item(5, 6).Value = "9"
$ws.Cells.Item(5, 7).Value = "3"
$ws.Cells.Item(5, 8).Value = "0"
$ws.Cells.Item(5, 9).Value = "133.33"
$ws.Cells.Item(5, 10).Value = "Mumbai Indians"
$ws.Cells.Item(5, 11).Value = "Sharjah"
$ws.Cells.Item(5, 12).Value = "October 05"
$ws.Cells.Item(5, 13).Value = "Mumbai won by 8 wickets (with 70 balls remaining)"

# Row 6
$ws.Cells.Item(6, 1).Value = "36th"
$ws.Cells.Item(6, 2).Value = "Rajasthan Royals"
$ws.Cells.Item(6, 3).Value = "Yashasvi Jaiswal"
$ws.Cells.Item(6, 4).Value = "c †Pant b Nortje"
$ws.Cells.Item(6, 5).Value = "5"
$ws.Cells.Item(6, 6).Value = "4"
$ws.Cells.Item(6, 7).Value = "0"
$ws.Cells.Item(6, 8).Value = "0"
$ws.Cells.Item(6, 9).Value = "125.00"
$ws.Cells.Item(6, 10).Value = "Delhi Capitals"
$ws.Cells.Item(6, 11).Value = "Abu Dhabi"
$ws.Cells.Item(6, 12).Value = "September 25"
$ws.Cells.Item(6, 13).Value = "Capitals won by 33 runs"

# Row 7
$ws.Cells.Item(7, 1).Value = "43rd"
$ws.Cells.Item(7, 2).Value = "Rajasthan Royals"
$ws.Cells.Item(7, 3).Value = "Yashasvi Jaiswal"
$ws.Cells.Item(7, 4).Value = "c Mohammed Siraj b Christian"
$ws.Cells.Item(7, 5).Value = "31"
$ws.Cells.Item(7, 6).Value = "22"
$ws.Cells.Item(7, 7).Value = "3"
$ws.Cells.Item(7, 8).Value = "2"
$ws.Cells.Item(7, 9).Value = "140.90"
$ws.Cells.Item(7, 10).Value = "Royal Challengers Bangalore"
$ws.Cells.Item(7, 11).Value = "Dubai (DSC)"
$ws.Cells.Item(7, 12).Value = "September 29"
$ws.Cells.Item(7, 13).Value = "RCB won by 7 wickets (with 17 balls remaining)"

# Row 8
$ws.Cells.Item(8, 1).Value = "18th"
$ws.Cells.Item(8, 2).Value = "Rajasthan Royals"
$ws.Cells.Item(8, 3).Value = "Yashasvi Jaiswal"
$ws.Cells.Item(8, 4).Value = "c sub (KL Nagarkoti) b Shivam Mavi"
$ws.Cells.Item(8, 5).Value = "22"
$ws.Cells.Item(8, 6).Value = "17"
$ws.Cells.Item(8, 7).Value = "5"
$ws.Cells.Item(8, 8).Value = "0"
$ws.Cells.Item(8, 9).Value = "129.41"
$ws.Cells.Item(8, 10).Value = "Kolkata Knight Riders"
$ws.Cells.Item(8, 11).Value = "Wankhede"
$ws.Cells.Item(8, 12).Value = "April 24"
$ws.Cells.Item(8, 13).Value = "Royals won by 6 wickets (with 7 balls remaining)"

# Row 9
$ws.Cells.Item(9, 1).Value = "32nd"
$ws.Cells.Item(9, 2).Value = "Rajasthan Royals"
$ws.Cells.Item(9, 3).Value = "Yashasvi Jaiswal"
$ws.Cells.Item(9, 4).Value = "c Agarwal b Harpreet Brar"
$ws.Cells.Item(9, 5).Value = "49"
$ws.Cells.Item(9, 6).Value = "36"
$ws.Cells.Item(9, 7).Value = "6"
$ws.Cells.Item(9, 8).Value = "2"
$ws.Cells.Item(9, 9).Value = "136.11"
$ws.Cells.Item(9, 10).Value = "Punjab Kings"
$ws.Cells.Item(9, 11).Value = "Dubai (DSC)"
$ws.Cells.Item(9, 12).Value = "September 21"
$ws.Cells.Item(9, 13).Value = "Royals won by 2 runs"

# Row 10
$ws.Cells.Item(10, 1).Value = "40th"
$ws.Cells.Item(10, 2).Value = "Rajasthan Royals"
$ws.Cells.Item(10, 3).Value = "Yashasvi Jaiswal"
$ws.Cells.Item(10, 4).Value = "b Sandeep Sharma"
$ws.Cells.Item(10, 5).Value = "36"
$ws.Cells.Item(10, 6).Value = "23"
$ws.Cells.Item(10, 7).Value = "5"
$ws.Cells.Item(10, 8).Value = "1"
$ws.Cells.Item(10, 9).Value = "156.52"
$ws.Cells.Item(10, 10).Value = "Sunrisers Hyderabad"
$ws.Cells.Item(10, 11).Value = "Dubai (DSC)"
$ws.Cells.Item(10, 12).Value = "September 27"
$ws.Cells.Item(10, 13).Value = "Sunrisers won by 7 wickets (with 9 balls remaining)"

# Row 11
$ws.Cells.Item(11, 1).Value = "54th"
$ws.Cells.Item(11, 2).Value = "Rajasthan Royals"
$ws.Cells.Item(11, 3).Value = "Yashasvi Jaiswal"
$ws.Cells.Item(11, 4).Value = "b Shakib Al Hasan"
$ws.Cells.Item(11, 5).Value = "0"
$ws.Cells.Item(11, 6).Value = "3"
$ws.Cells.Item(11, 7).Value = "0"
$ws.Cells.Item(11, 8).Value = "0"
$ws.Cells.Item(11, 9).Value = "0.00"
$ws.Cells.Item(11, 10).Value = "Kolkata Knight Riders"
$ws.Cells.Item(11, 11).Value = "Sharjah"
$ws.Cells.Item(11, 12).Value = "October 07"
$ws.Cells.Item(11, 13).Value = "KKR won by 86 runs"

# Drop back to the default "Normal" style so no explicit number-format
# style index lingers on the cells (matches the original, style-less cells).
$target.Style = "Normal"

